$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed rows 142/143 with the formatting (font/fill/border/number-format)
# --- of the last existing data row (141), then overwrite the values below.
$ws.Range("A141:M141").Copy()
$ws.Range("A142:M142").PasteSpecial(-4122)
$ws.Range("A141:M141").Copy()
$ws.Range("A143:M143").PasteSpecial(-4122)

# ===================== Row 142 =====================
$ws.Cells.Item(142, 1).Value = "$([char]9679)"
# D142 must stay text ("10323"), but its style is General-formatted, so a
# plain numeric-looking string would be coerced to a number - force text
# entry with a leading apostrophe, then restore the row's normal style.
$ws.Cells.Item(142, 4).Value = "'10323"
$ws.Range("D141").Copy()
$ws.Range("D142").PasteSpecial(-4122)
$ws.Cells.Item(142, 5).Value = "Location:Country"
$ws.Cells.Item(142, 6).Value = "2: 3521"
$ws.Cells.Item(142, 7).Value = "2: 3525"
$ws.Cells.Item(142, 8).Value = 0
$ws.Cells.Item(142, 9).Value = "Japan"
$ws.Cells.Item(142, 10).Value = 5
$ws.Cells.Item(142, 11).Value = 0.037268932617769826
$ws.Cells.Item(142, 12).Value = "dattaray"
$ws.Cells.Item(142, 13).Value = "9/17/2019 12:14:42"
$ws.Rows.Item(142).RowHeight = 15.6

# ===================== Row 143 =====================
$ws.Cells.Item(143, 1).Value = "$([char]9679)"
$ws.Cells.Item(143, 4).Value = "'10323"
$ws.Range("D141").Copy()
$ws.Range("D143").PasteSpecial(-4122)
$ws.Cells.Item(143, 5).Value = "Exclusion: Not first case in country"
$ws.Cells.Item(143, 6).Value = "3: 424"
$ws.Cells.Item(143, 7).Value = "3: 497"
$ws.Cells.Item(143, 8).Value = 0
$ws.Cells.Item(143, 9).Value = "To our knowledge, this is `nthe first report of a fatality caused by ISMRK."
$ws.Cells.Item(143, 10).Value = 74
$ws.Cells.Item(143, 11).Value = 0.55158020274299346
$ws.Cells.Item(143, 12).Value = "dattaray"
$ws.Cells.Item(143, 13).Value = "9/17/2019 12:15:52"
$ws.Rows.Item(143).RowHeight = 25.8

# Column I got wider to fit the new, longer comment text.
$ws.Columns.Item(9).ColumnWidth = 50.8
